$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.790.22"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "2.818.86"
$ws.Range("E3").Value = "  +1.55%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.95"
$ws.Range("E5").Value = "  +3.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.80"
$ws.Range("E6").Value = "  -3.01%  "

$ws.Range("E7").Value = "  +3.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  +3.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.77"
$ws.Range("E10").Value = "  -4.23%  "

$ws.Range("E11").Value = "  +0.03%  "

$ws.Range("E12").Value = "  +0.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.91"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("E14").Value = "  +1.39%  "

$ws.Range("D15").Value = "3.260.25"
$ws.Range("E15").Value = "  +1.67%  "

$ws.Range("D16").Value = "2.812.45"
$ws.Range("E16").Value = "  +1.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.920"
$ws.Range("E17").Value = "  +4.51%  "

$ws.Range("D18").Value = "51.722.00"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.58"
$ws.Range("E19").Value = "  +7.98%  "

$ws.Range("E20").Value = "  -2.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.35"
$ws.Range("E21").Value = "  +0.85%  "

$ws.Range("D22").Value = "0.0₃0991"
$ws.Range("E22").Value = "  +1.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.01"
$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.86"
$ws.Range("E24").Value = "  -1.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +1.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.91"
$ws.Range("E26").Value = "  +1.55%  "

$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.30"
$ws.Range("E28").Value = "  +0.92%  "

$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0491"
$ws.Range("E30").Value = "  +28.72%  "

$ws.Range("E31").Value = "  +0.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.57"
$ws.Range("E32").Value = "  +4.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.50"
$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.86"
$ws.Range("E34").Value = "  +2.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.39"
$ws.Range("E35").Value = "  +9.29%  "

$ws.Range("E36").Value = "  +3.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.29"
$ws.Range("E38").Value = "  +2.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.03"
$ws.Range("E39").Value = "  -3.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.32"
$ws.Range("E40").Value = "  -3.35%  "

$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("E42").Value = "  -4.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.23"
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "124.72"
$ws.Range("E44").Value = "  -1.96%  "

$ws.Range("E45").Value = "  -2.44%  "

$ws.Range("D46").Value = "2.093.53"
$ws.Range("E46").Value = "  +1.22%  "

$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("E48").Value = "  +1.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.98"
$ws.Range("E49").Value = "  +8.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.988"
$ws.Range("E50").Value = "  +10.21%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "61.08"
$ws.Range("E51").Value = "  +2.77%  "

